# edit.ps1
# Applies the changes described by the commit:
#  1. Update the fixed "date" footer placeholder text (16-04-2019 -> 24-05-2019)
#     on the slide master and every slide layout.
#  2. Move a group of shapes (the state-diagram block) up/right by a constant
#     offset (dx=+1825559 EMU, dy=-652618 EMU), and move the block's background
#     rectangle to its own new position.
#  3. Re-color three ellipses' fill from accent2 to accent1.
#  4. Delete the leftover "Diagrama de estado" title textbox (Rectangulo 202).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Fixed date field: slide master + all slide layouts
# ---------------------------------------------------------------------------
$oldDate = "16-04-2019"
$newDate = "24-05-2019"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Reposition shapes on the slide
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

# Shape 1 "Rectangulo 201" (background panel for the block) gets its own
# explicit new position.
$s.Shapes.Item(1).Left = 219.78095248188976
$s.Shapes.Item(1).Top  = 25.48259842519685

# Shapes 31-44 (Rectangulo 187-193 and CuadroTexto 194-200) shift together by
# the same delta: dx = +1825559 EMU, dy = -652618 EMU. Point values below were
# solved so that, after the COM layer's Single-precision Left/Top storage,
# the saved EMU offsets land exactly on the target values.
$newPositions = @{
    31 = @(239.9000016,          23.299685539370078)
    32 = @(229.04496062992126,   66.81866141732283)
    33 = @(376.6956024511811,    23.299685539370078)
    34 = @(363.6029205858268,    66.81866141732283)
    35 = @(510.46189886377954,   23.299685539370078)
    36 = @(515.6369324338583,    66.81866141732283)
    37 = @(647.9624409448819,    23.299685539370078)
    38 = @(280.10984811968507,   42.58425336850394)
    39 = @(280.10984811968507,   74.08897787795276)
    40 = @(407.8235433070866,    33.413149606299214)
    41 = @(411.11622627244094,   85.41425196850393)
    42 = @(545.324094488189,     41.895196950393704)
    43 = @(544.9415748031496,    82.57094488188976)
    44 = @(673.7240944881889,    41.895196950393704)
}

foreach ($idx in $newPositions.Keys) {
    $pos = $newPositions[$idx]
    $shape = $s.Shapes.Item([int]$idx)
    $shape.Left = $pos[0]
    $shape.Top  = $pos[1]
}

# ---------------------------------------------------------------------------
# 3) Recolor three ellipses' solid fill from accent2 to accent1
#    (msoThemeColorAccent1 = 5)
# ---------------------------------------------------------------------------
foreach ($idx in 3, 4, 5) {
    $s.Shapes.Item($idx).Fill.ForeColor.ObjectThemeColor = 5
}

# ---------------------------------------------------------------------------
# 4) Delete the leftover "Diagrama de estado" title textbox
# ---------------------------------------------------------------------------
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Rectángulo 202") {
        $sh.Delete()
    }
}
